$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of row number -> new values for columns B, C, E, F, G, J, K, L, M, O
# (columns D, H, I, N remain 0 and are unchanged)
$cols = @("B","C","E","F","G","J","K","L","M","O")

$data = @{
    2 = @(13.3167191879516, 11.1427967298722, 15.7042566683154, 42.60373632216616, 3.69459859876229, 9.136431601577488, 9.093964355766161, 11.82890564470812, 15.7575500699376, 29.28222013225299)
    3 = @(13.12966339586001, 11.15864958788797, 15.72009619313027, 42.66176209332836, 3.696424285542953, 9.133706931412892, 8.955846068797781, 11.82613124004907, 15.72490543507519, 29.37695315391028)
    4 = @(13.01582749366839, 11.16895668298351, 15.73179125538978, 42.70593196365606, 3.697604829751401, 9.132098849912024, 8.871333712402505, 11.82583430877352, 15.70691683953704, 29.44015539324957)
    5 = @(12.96974939383114, 11.17330149154136, 15.73705307243111, 42.72607795016197, 3.698100935718328, 9.131460015520929, 8.837008895105464, 11.82606838189556, 15.70010838364225, 29.46717615570679)
    6 = @(12.96211854942491, 11.17403168809694, 15.73795676934194, 42.72955275791765, 3.698184222553177, 9.131354937032256, 8.831317412274805, 11.82612873354823, 15.69900952033449, 29.47173932782939)
    7 = @(13.01520473651435, 11.16901469263702, 15.73186020895416, 42.70619497120096, 3.697611459514436, 9.132090167468862, 8.870870277358565, 11.82583602612861, 15.70682289796011, 29.44051468268834)
    8 = @(13.25204440700425, 11.14814403361531, 15.70930981713462, 42.62197002044943, 3.695215760843749, 9.135478764092811, 9.046304926001403, 11.82765795597268, 15.74587083713301, 29.31383833805593)
    9 = @(13.72213212383937, 11.11174801429889, 15.68068222540411, 42.52462685792025, 3.690988303106466, 9.142632000154711, 9.390871102247175, 11.84232738573347, 15.83851015402305, 29.10542494459252)
    10 = @(14.06772966206837, 11.08774506707434, 15.66910813827512, 42.49448674937508, 3.688166245565132, 9.148189042407633, 9.642008825054733, 11.85977962941452, 15.91602428759903, 28.97674404833002)
    11 = @(14.22433484447649, 11.0774145588915, 15.66588498905964, 42.4897536065972, 3.686943422652241, 9.150780824904601, 9.755342732300123, 11.86914701896629, 15.95326027596865, 28.92352252673362)
    12 = @(14.28349592569703, 11.0735868961088, 15.66495695643326, 42.48925043638459, 3.686489087461226, 9.151771286829073, 9.798090342709949, 11.87289755672784, 15.96763721542397, 28.90413436185179)
    13 = @(14.27076168576435, 11.07440750918751, 15.66514383453969, 42.48930150372552, 3.686586549464668, 9.151557575157629, 9.788892016396627, 11.87208080326397, 15.96452870941763, 28.90827587747854)
    14 = @(14.22920521112368, 11.07709796777142, 15.66580278380368, 42.48968638664499, 3.686905869676845, 9.150862131217387, 9.758863216405601, 11.86945151560348, 15.95443757805147, 28.92191209865773)
    15 = @(14.20373060792503, 11.07875691611545, 15.66624446698517, 42.49008995985731, 3.687102597041983, 9.150437319977003, 9.740446517499068, 11.86786741453065, 15.94829224853581, 28.93036442659178)
    16 = @(14.05747820748116, 11.08843200266681, 15.66935977371542, 42.49497658752299, 3.688247382647546, 9.148020933572061, 9.634580540607733, 11.85919601441087, 15.91362992936831, 28.98032927832874)
    17 = @(13.96756193782163, 11.09451784208946, 15.67179313373342, 42.50027272601679, 3.688965250688545, 9.14655479163342, 9.569374177858194, 11.85424076831354, 15.89286621734918, 29.01234355096637)
    18 = @(13.91579031688879, 11.09807367438063, 15.67338500997194, 42.50416410720302, 3.689383888371002, 9.145717526369037, 9.531785778324418, 11.85152520937236, 15.88110978582725, 29.03125761542271)
    19 = @(13.89825371770446, 11.09928714720873, 15.67395704774517, 42.50562686650264, 3.68952661887329, 9.145435082049316, 9.519045865900292, 11.85062894117602, 15.87716148022232, 29.03774748206092)
    20 = @(13.97713968290908, 11.0938642610959, 15.67151420688334, 42.49962148069182, 3.688888238756505, 9.146710243057335, 9.576324425308437, 11.85475435014263, 15.89505732087201, 29.00888378864633)
    21 = @(14.24141563053552, 11.07630543044844, 15.6656013053689, 42.48953836766331, 3.686811841234087, 9.151066156806287, 9.767688302797135, 11.87021830042887, 15.95739414515294, 28.91788602026351)
    22 = @(14.41328421890763, 11.06532078510296, 15.66344123579262, 42.49046136921805, 3.68550561037993, 9.153965435522258, 9.891750990146507, 11.88150889304762, 15.99974273699952, 28.86287717613535)
    23 = @(14.32165007039729, 11.07113868116685, 15.6644385518767, 42.48928210512142, 3.686198134762527, 9.152413295108753, 9.825640686220977, 11.875375269268, 15.97699590957498, 28.89182762811116)
    24 = @(13.97280982177633, 11.09415956748329, 15.67163970875612, 42.49991327159117, 3.688923037389327, 9.146639945884411, 9.573182528179007, 11.85452174450011, 15.89406615891533, 29.01044636289932)
    25 = @(13.59469171762171, 11.12111161789503, 15.68676194954521, 42.54369403817824, 3.692081880864829, 9.140643302001783, 9.297844516610176, 11.83717980251333, 15.81176276332315, 29.15751827456683)
}

foreach ($row in $data.Keys) {
    $values = $data[$row]
    for ($i = 0; $i -lt $cols.Length; $i++) {
        $cellRef = "{0}{1}" -f $cols[$i], $row
        $ws.Range($cellRef).Value = $values[$i]
    }
}

Write-Host "Updated loading_percent values for rows 2-25"
